# Apply the Saldo.xlsx update:
#  - Insert a new row for account 004504449 / KELMA / 1060.21 right above the
#    004415557 / FILIPE row (keeping the sheet sorted by descending Saldo).
#  - Remove the old 004504449 / KELMA / 60.21 row further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row where account 004415557 currently lives; insert the new
# row directly above it, shifting that row (and everything below) down.
$insertRow = $ws.Cells.Find("004415557").Row
$ws.Rows.Item($insertRow).Insert()

# Write the "Conta" value as text (not a number) so the leading zeros are
# preserved, same as every other account-number cell in the sheet.
$ws.Cells.Item($insertRow, 1).Formula = "=""004504449"""
$ws.Cells.Item($insertRow, 1).Copy()
$ws.Cells.Item($insertRow, 1).PasteSpecial(-4163)

$ws.Cells.Item($insertRow, 2).Value = "KELMA"
$ws.Cells.Item($insertRow, 3).Value = 1060.21

# Now locate (and remove) the stale 004504449 / 60.21 row, which has
# shifted one row further down because of the insert above.
$oldRow = $ws.Cells.Find("004504449", $ws.Cells.Item($insertRow, 1)).Row
$ws.Rows.Item($oldRow).Delete()
